# Append a new row (row 6) to Sheet1, matching a new submission captured at
# 01/05/2025 04:37:10 PM: احمد / 22 / الجزائري / الرحلة 1 / C3 / NRC.
#
# The sheet stores every cell as text (SheeJS-style export with a
# numberStoredAsText ignoredError suppression over the used range), so every
# value written below must land as a text cell - even the ones that look
# like plain numbers ("22") or the blank notes cell. A leading "'" forces
# Excel's quote-prefix (text) interpretation for values that would otherwise
# be parsed as a number; re-applying the "Normal" style afterwards strips the
# transient quote-prefix cell style back off so the cell matches its
# neighbours (no style override), while the stored value/type stay text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "'"
$ws.Range("A6").Style = "Normal"

$ws.Range("B6").Value = "احمد"

$ws.Range("C6").Value = "'22"
$ws.Range("C6").Style = "Normal"

$ws.Range("D6").Value = "الجزائري"
$ws.Range("E6").Value = "الرحلة 1"
$ws.Range("F6").Value = "C3"
$ws.Range("G6").Value = "NRC"
$ws.Range("H6").Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٤:٣٧:١٠ م"
